# Apply cell updates from the "Updated symbol list" GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "305.01"
Set-TextValue "E2" "-0.29%"
Set-TextValue "G2" "23"

# Row 3
Set-TextValue "D3" "35.90"
Set-TextValue "E3" "-1.07%"
Set-TextValue "G3" "23"

# Row 4
Set-TextValue "D4" "4.998"
Set-TextValue "E4" "-2.16%"
Set-TextValue "G4" "23"

# Row 5
Set-TextValue "D5" "0.08062"
Set-TextValue "E5" "-0.20%"
Set-TextValue "G5" "23"

# Row 6
Set-TextValue "D6" "1.898"
Set-TextValue "E6" "-1.95%"
Set-TextValue "G6" "23"

# Row 7
Set-TextValue "D7" "7.852"
Set-TextValue "E7" "1.56%"
Set-TextValue "G7" "23"

# Row 8
Set-TextValue "D8" "0.9342"
Set-TextValue "E8" "0.30%"
Set-TextValue "G8" "23"

# Row 9
Set-TextValue "D9" "0.1291"
Set-TextValue "E9" "-11.57%"
Set-TextValue "G9" "23"

# Row 10
Set-TextValue "D10" "0.1900"
Set-TextValue "G10" "23"

# Row 11
Set-TextValue "D11" "0.09191"
Set-TextValue "E11" "0.97%"
Set-TextValue "G11" "23"

# Row 12
Set-TextValue "D12" "0.03508"
Set-TextValue "E12" "-1.18%"
Set-TextValue "G12" "23"

# Row 13
Set-TextValue "D13" "0.09894"
Set-TextValue "E13" "1.03%"
Set-TextValue "G13" "23"

# Row 14
Set-TextValue "D14" "0.001426"
Set-TextValue "E14" "-0.49%"
Set-TextValue "G14" "23"

# Row 15
Set-TextValue "D15" "0.006548"
Set-TextValue "E15" "12.80%"
Set-TextValue "G15" "23"

# Row 16
Set-TextValue "D16" "3.612"
Set-TextValue "E16" "2.54%"
Set-TextValue "G16" "23"

# Row 17
Set-TextValue "D17" "4.150"
Set-TextValue "E17" "0.78%"
Set-TextValue "G17" "23"

# Row 18
Set-TextValue "D18" "3.107"
Set-TextValue "E18" "5.98%"
Set-TextValue "G18" "23"

# Row 19
Set-TextValue "E19" "0.55%"
Set-TextValue "G19" "23"

# Row 20
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D20" "0.1335"
Set-TextValue "E20" "2.57%"
Set-TextValue "G20" "23"

# Row 21
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D21" "5.236"
Set-TextValue "E21" "3.77%"
Set-TextValue "G21" "23"

# Row 22
Set-TextValue "G22" "23"

# Row 23
Set-TextValue "D23" "0.04419"
Set-TextValue "E23" "-2.51%"
Set-TextValue "G23" "23"

# Row 24
Set-TextValue "E24" "1.87%"
Set-TextValue "G24" "23"

# Row 25
Set-TextValue "D25" "0.004716"
Set-TextValue "E25" "-2.78%"
Set-TextValue "G25" "23"

# Row 26
Set-TextValue "E26" "4.25%"
Set-TextValue "G26" "23"

# Row 27
Set-TextValue "E27" "-29.73%"
Set-TextValue "G27" "23"

# Row 28
Set-TextValue "G28" "23"

# Row 29
Set-TextValue "G29" "23"

# Row 30
Set-TextValue "G30" "23"

# Row 31
Set-TextValue "G31" "23"

# Row 32
Set-TextValue "G32" "23"

# Row 33
Set-TextValue "G33" "23"

# Row 34
Set-TextValue "G34" "23"

# Row 35
Set-TextValue "G35" "23"

# Row 36
Set-TextValue "G36" "23"

# Row 37
Set-TextValue "G37" "23"

# Row 38
Set-TextValue "G38" "23"

# Row 39
Set-TextValue "D39" "0.01954"
Set-TextValue "E39" "-0.90%"
Set-TextValue "G39" "23"

# Row 40
Set-TextValue "D40" "0.05171"
Set-TextValue "E40" "7.07%"
Set-TextValue "G40" "23"

# Row 41
Set-TextValue "D41" "0.007552"
Set-TextValue "E41" "0.60%"
Set-TextValue "G41" "23"

# Row 42
Set-TextValue "D42" "0.01020"
Set-TextValue "E42" "-8.57%"
Set-TextValue "G42" "23"

# Row 43
Set-TextValue "D43" "0.1373"
Set-TextValue "E43" "0.77%"
Set-TextValue "G43" "23"

# Row 44
Set-TextValue "E44" "9.11%"
Set-TextValue "G44" "23"

# Row 45
Set-TextValue "D45" "0.01072"
Set-TextValue "E45" "8.58%"
Set-TextValue "G45" "23"

# Row 46
Set-TextValue "D46" "0.00006360"
Set-TextValue "E46" "-0.56%"
Set-TextValue "G46" "23"

# Row 47
Set-TextValue "E47" "-0.17%"
Set-TextValue "G47" "23"

# Row 48
Set-TextValue "D48" "65.22"
Set-TextValue "E48" "0.85%"
Set-TextValue "G48" "23"

# Row 49
Set-TextValue "E49" "39.12%"
Set-TextValue "G49" "23"

# Row 50
Set-TextValue "E50" "-0.17%"
Set-TextValue "G50" "23"

# Row 51
Set-TextValue "E51" "-0.17%"
Set-TextValue "G51" "23"
